# Daily attendance processing - 2025-12-06 09:48:29
# Reorders the comma-separated "Recorded By" values in column G
# (reverses the order of the names/emails listed in each cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$rows = @(2,4,5,7,8,11,17,28,30,31,33,34,37,43,54,56,57,59,60,63,69,80,81,82,93,94,96,106,107,108,119,120,122,132,133,134,145,146,148)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = Recorded By
    $current = [string]$cell.Value2
    $parts = $current -split ',\s*'
    $reversed = $parts[($parts.Length - 1)..0]
    $cell.Value2 = [string]::Join(', ', $reversed)
}
